$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Row 2 values (Beta) ---
$ws.Range("C2").Value = 19.12075701903682
$ws.Range("D2").Value = 0.006944391349456487
$ws.Range("E2").Value = 0.01982943797740053
$ws.Range("F2").Value = 9.824617691558897
$ws.Range("G2").Value = 9.621279507202052
$ws.Range("H2").Value = 10.03041475721893
$ws.Range("I2").Value = 0.002527942562296337
$ws.Range("J2").Value = 0.002253059092366989
$ws.Range("K2").Value = 0.002865412155953543
$ws.Range("L2").Value = 0.008890111058752601
$ws.Range("M2").Value = 0.008696115606330237
$ws.Range("N2").Value = 0.009092402267896798

# --- Update existing Row 3 values (Gamma) ---
$ws.Range("C3").Value = 0.04981522627320694
$ws.Range("D3").Value = 0.04815098319456564
$ws.Range("E3").Value = 0.0499839736740351
$ws.Range("F3").Value = 0.0466430919665593
$ws.Range("G3").Value = 0.04636586929448727
$ws.Range("H3").Value = 0.04694978669319402
$ws.Range("I3").Value = 0.0450833632713068
$ws.Range("J3").Value = 0.04481643278537908
$ws.Range("K3").Value = 0.04537777210227339
$ws.Range("L3").Value = 0.04669361636918469
$ws.Range("M3").Value = 0.04641639815915267
$ws.Range("N3").Value = 0.04700027433807879

# --- Add new Row 4 (Beta + Gamma) ---
# Copy formatting (border/font/alignment) from A2 (style index 1) onto A4
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2

$ws.Range("B4").Value = "Beta + Gamma"

$ws.Range("C4").Value = 19.17057224531003
$ws.Range("D4").Value = 0.05509537454402212
$ws.Range("E4").Value = 0.06981341165143562
$ws.Range("F4").Value = 9.871260783525457
$ws.Range("G4").Value = 9.667645376496539
$ws.Range("H4").Value = 10.07736454391212
$ws.Range("I4").Value = 0.04761130583360314
$ws.Range("J4").Value = 0.04706949187774607
$ws.Range("K4").Value = 0.04824318425822693
$ws.Range("L4").Value = 0.0555837274279373
$ws.Range("M4").Value = 0.0551125137654829
$ws.Range("N4").Value = 0.05609267660597559
